# Weekly Acelga price update for "Terminal La Palmera de La Serena".
# A new week of data (2 rows: Primera / Segunda quality) is inserted
# right above the existing row 212, pushing all subsequent rows down by
# two (old row 212 -> new row 214, ... old row 357 -> new row 359).
#
# The new rows reuse the price-tier figures (K/L/M/N/O/P/Q, etc.) of the
# week that used to sit at rows 212-213, only the date (D) and the
# volume (J) differ for the newly reported week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new week: insert two blank rows at 212:213, shifting
# the rest of the table (212:357 -> 214:359) down.
$ws.Rows("212:213").Insert()

# Seed the two new rows with the same row layout/styling/data as the
# row that is now directly below them (i.e. the old row 212/213 data,
# now living at 214/215) and then overwrite the date + volume with the
# values for the newly reported week.
$ws.Range("A214:R215").Copy()
$ws.Range("A212").PasteSpecial()

$ws.Range("D212").Value = 44634
$ws.Range("J212").Value = 2400

$ws.Range("D213").Value = 44634
$ws.Range("J213").Value = 1300
